$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Terminal Hortofrutícola
# Agro Chillán / Zapallo italiano. Insert a new row at 306 (pushing the
# existing rows 306-335 down to 307-336) and populate it with the new
# record, matching the rest of the subset's columns.
$ws.Rows("306").Insert()

$ws.Cells.Item(306, 1).Value = 7
$ws.Cells.Item(306, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(306, 3).Value = "Ñuble"
$ws.Cells.Item(306, 4).Value = 45132
$ws.Cells.Item(306, 5).Value = 16
$ws.Cells.Item(306, 6).Value = 100112032
$ws.Cells.Item(306, 7).Value = "Zapallo italiano"
$ws.Cells.Item(306, 8).Value = "Sin especificar"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 80
$ws.Cells.Item(306, 11).Value = 14000
$ws.Cells.Item(306, 12).Value = 14000
$ws.Cells.Item(306, 13).Value = 14000
$ws.Cells.Item(306, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(306, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(306, 16).Value = 280
$ws.Cells.Item(306, 17).Value = 50
$ws.Cells.Item(306, 18).Value = "Hortaliza"
